$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DAMSLTag (column I) and DialogAct (column J) values
# following a re-run of SGNN dialog act annotation.

$ws.Range("I5").Value = "%"
$ws.Range("J5").Value = "Uninterpretable"
$ws.Range("I23").Value = "sv"
$ws.Range("J23").Value = "Statement-opinion"
$ws.Range("I34").Value = "b"
$ws.Range("J34").Value = "Acknowledge (Backchannel)"
$ws.Range("I52").Value = "aa"
$ws.Range("J52").Value = "Agree/Accept"
$ws.Range("I58").Value = "aa"
$ws.Range("J58").Value = "Agree/Accept"
$ws.Range("I66").Value = "aa"
$ws.Range("J66").Value = "Agree/Accept"
$ws.Range("I83").Value = "b"
$ws.Range("J83").Value = "Acknowledge (Backchannel)"
$ws.Range("I86").Value = "sd"
$ws.Range("J86").Value = "Statement-non-opinion"
$ws.Range("I93").Value = "sv"
$ws.Range("J93").Value = "Statement-opinion"
$ws.Range("I111").Value = "sd"
$ws.Range("J111").Value = "Statement-non-opinion"
$ws.Range("I117").Value = "b"
$ws.Range("J117").Value = "Acknowledge (Backchannel)"
$ws.Range("I137").Value = "sv"
$ws.Range("J137").Value = "Statement-opinion"
$ws.Range("I139").Value = "b"
$ws.Range("J139").Value = "Acknowledge (Backchannel)"
$ws.Range("I154").Value = "%"
$ws.Range("J154").Value = "Uninterpretable"
$ws.Range("I162").Value = "ba"
$ws.Range("J162").Value = "Appreciation"
$ws.Range("I163").Value = "sd"
$ws.Range("J163").Value = "Statement-non-opinion"
$ws.Range("I164").Value = "ba"
$ws.Range("J164").Value = "Appreciation"
$ws.Range("I165").Value = "ba"
$ws.Range("J165").Value = "Appreciation"
$ws.Range("I172").Value = "sd"
$ws.Range("J172").Value = "Statement-non-opinion"
$ws.Range("I194").Value = "sv"
$ws.Range("J194").Value = "Statement-opinion"
$ws.Range("I219").Value = "sd"
$ws.Range("J219").Value = "Statement-non-opinion"
$ws.Range("I245").Value = "sd"
$ws.Range("J245").Value = "Statement-non-opinion"
$ws.Range("I257").Value = "b"
$ws.Range("J257").Value = "Acknowledge (Backchannel)"
$ws.Range("I266").Value = "sd"
$ws.Range("J266").Value = "Statement-non-opinion"
$ws.Range("I271").Value = "aa"
$ws.Range("J271").Value = "Agree/Accept"
$ws.Range("I310").Value = "sv"
$ws.Range("J310").Value = "Statement-opinion"
$ws.Range("I340").Value = "aa"
$ws.Range("J340").Value = "Agree/Accept"
$ws.Range("I341").Value = "sd"
$ws.Range("J341").Value = "Statement-non-opinion"
$ws.Range("I363").Value = "%"
$ws.Range("J363").Value = "Uninterpretable"
$ws.Range("I365").Value = "b"
$ws.Range("J365").Value = "Acknowledge (Backchannel)"
$ws.Range("I368").Value = "ba"
$ws.Range("J368").Value = "Appreciation"
$ws.Range("I376").Value = "sv"
$ws.Range("J376").Value = "Statement-opinion"
$ws.Range("I387").Value = "ba"
$ws.Range("J387").Value = "Appreciation"
$ws.Range("I390").Value = "b"
$ws.Range("J390").Value = "Acknowledge (Backchannel)"
$ws.Range("I398").Value = "aa"
$ws.Range("J398").Value = "Agree/Accept"
$ws.Range("I411").Value = "aa"
$ws.Range("J411").Value = "Agree/Accept"
$ws.Range("I415").Value = "b"
$ws.Range("J415").Value = "Acknowledge (Backchannel)"
$ws.Range("I421").Value = "sv"
$ws.Range("J421").Value = "Statement-opinion"
$ws.Range("I443").Value = "sd"
$ws.Range("J443").Value = "Statement-non-opinion"
$ws.Range("I444").Value = "sv"
$ws.Range("J444").Value = "Statement-opinion"
$ws.Range("I449").Value = "qy"
$ws.Range("J449").Value = "Yes-No-Question"
$ws.Range("I451").Value = "sv"
$ws.Range("J451").Value = "Statement-opinion"
$ws.Range("I459").Value = "sv"
$ws.Range("J459").Value = "Statement-opinion"
$ws.Range("I460").Value = "aa"
$ws.Range("J460").Value = "Agree/Accept"
$ws.Range("I461").Value = "b"
$ws.Range("J461").Value = "Acknowledge (Backchannel)"
$ws.Range("I462").Value = "sd"
$ws.Range("J462").Value = "Statement-non-opinion"
$ws.Range("I463").Value = "sv"
$ws.Range("J463").Value = "Statement-opinion"
$ws.Range("I472").Value = "sv"
$ws.Range("J472").Value = "Statement-opinion"
$ws.Range("I480").Value = "b"
$ws.Range("J480").Value = "Acknowledge (Backchannel)"
$ws.Range("I503").Value = "aa"
$ws.Range("J503").Value = "Agree/Accept"
$ws.Range("I506").Value = "aa"
$ws.Range("J506").Value = "Agree/Accept"
$ws.Range("I508").Value = "aa"
$ws.Range("J508").Value = "Agree/Accept"
$ws.Range("I509").Value = "aa"
$ws.Range("J509").Value = "Agree/Accept"
$ws.Range("I510").Value = "aa"
$ws.Range("J510").Value = "Agree/Accept"
$ws.Range("I516").Value = "sv"
$ws.Range("J516").Value = "Statement-opinion"
$ws.Range("I536").Value = "aa"
$ws.Range("J536").Value = "Agree/Accept"
$ws.Range("I537").Value = "b"
$ws.Range("J537").Value = "Acknowledge (Backchannel)"
$ws.Range("I543").Value = "sv"
$ws.Range("J543").Value = "Statement-opinion"
$ws.Range("I556").Value = "sv"
$ws.Range("J556").Value = "Statement-opinion"
$ws.Range("I559").Value = "b"
$ws.Range("J559").Value = "Acknowledge (Backchannel)"
$ws.Range("I566").Value = "ba"
$ws.Range("J566").Value = "Appreciation"
$ws.Range("I568").Value = "sd"
$ws.Range("J568").Value = "Statement-non-opinion"
$ws.Range("I572").Value = "aa"
$ws.Range("J572").Value = "Agree/Accept"
$ws.Range("I574").Value = "sv"
$ws.Range("J574").Value = "Statement-opinion"
$ws.Range("I576").Value = "sv"
$ws.Range("J576").Value = "Statement-opinion"
$ws.Range("I578").Value = "ba"
$ws.Range("J578").Value = "Appreciation"
$ws.Range("I579").Value = "aa"
$ws.Range("J579").Value = "Agree/Accept"
$ws.Range("I581").Value = "ba"
$ws.Range("J581").Value = "Appreciation"
$ws.Range("I596").Value = "b"
$ws.Range("J596").Value = "Acknowledge (Backchannel)"
$ws.Range("I604").Value = "sd"
$ws.Range("J604").Value = "Statement-non-opinion"
$ws.Range("I605").Value = "aa"
$ws.Range("J605").Value = "Agree/Accept"
$ws.Range("I614").Value = "sd"
$ws.Range("J614").Value = "Statement-non-opinion"
$ws.Range("I618").Value = "%"
$ws.Range("J618").Value = "Uninterpretable"
$ws.Range("I620").Value = "b"
$ws.Range("J620").Value = "Acknowledge (Backchannel)"
